$d = $word.ActiveDocument

$pairs = @(
    @("825×9=", "422×6="),
    @("301×7=", "153×6="),
    @("686×3=", "353×4="),
    @("390×8=", "540×7="),
    @("613×3=", "695×6="),
    @("941×3=", "177×8="),
    @("295×6=", "159×2="),
    @("966×5=", "286×3="),
    @("170×3=", "391×9="),
    @("611×5=", "735×2="),
    @("161×6=", "148×4="),
    @("289×8=", "274×4="),
    @("471×6=", "957×6="),
    @("101×8=", "191×6="),
    @("672×2=", "606×3="),
    @("743×7=", "403×8="),
    @("185×2=", "189×4="),
    @("351×7=", "701×8="),
    @("357×3=", "808×7="),
    @("737×2=", "812×8="),
    @("925×2=", "348×2="),
    @("443×2=", "963×3="),
    @("305×9=", "195×3="),
    @("766×7=", "337×4="),
    @("750×8=", "287×6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
